$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 11371486
$ws.Range("I51").Value = 41671170
$ws.Range("J51").Value = 9105.875
$ws.Range("K51").Value = 41671170
$ws.Range("L51").Value = 9105.875
$ws.Range("M51").Value = -41670686
$ws.Range("N51").Value = -10073.875
$ws.Range("H81").Value = 60000
$ws.Range("J81").Value = 60000
$ws.Range("L81").Value = 60000
$ws.Range("N81").Value = -61996
$ws.Range("H84").Value = 60000
$ws.Range("J84").Value = 60000
$ws.Range("L84").Value = 180000
$ws.Range("N84").Value = -189984
$ws.Range("H98").Value = 2696.4138
$ws.Range("I98").Value = 2278.6086
$ws.Range("K98").Value = 2278.6086
$ws.Range("M98").Value = -780.6086
$ws.Range("H103").Value = 992.16
$ws.Range("J103").Value = 1071.7
$ws.Range("L103").Value = 3215.1
$ws.Range("N103").Value = -4387.1
$ws.Range("H113").Value = 4155.294
$ws.Range("I113").Value = 3765.75
$ws.Range("K113").Value = 3765.75
$ws.Range("M113").Value = -511.75
$ws.Range("H122").Value = 2696.4138
$ws.Range("I122").Value = 2278.6086
$ws.Range("K122").Value = 6835.825800000001
$ws.Range("M122").Value = -4385.825800000001
$ws.Range("H132").Value = 2872660
$ws.Range("J132").Value = 1797
$ws.Range("L132").Value = 5391
$ws.Range("N132").Value = -10451
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21136.908
$ws.Range("I32").Value = 21136.908
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 21136.908
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -20849.908
$ws.Range("H69").Value = 350000
$ws.Range("J69").Value = 350000
$ws.Range("L69").Value = 350000
$ws.Range("N69").Value = -351498
$ws.Range("H72").Value = 350000
$ws.Range("J72").Value = 350000
$ws.Range("L72").Value = 1050000
$ws.Range("N72").Value = -1057488
$ws.Range("H74").Value = 359872.25
$ws.Range("I74").Value = 376114.25
$ws.Range("K74").Value = 376114.25
$ws.Range("M74").Value = -375240.25
$ws.Range("H77").Value = 359872.25
$ws.Range("I77").Value = 376114.25
$ws.Range("K77").Value = 1880571.25
$ws.Range("M77").Value = -1876203.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2291.6667
$ws.Range("I54").Value = 2291.6667
$ws.Range("K54").Value = 2291.6667
$ws.Range("M54").Value = -1807.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 18408
$ws.Range("I19").Value = 112
$ws.Range("K19").Value = 112
$ws.Range("M19").Value = 58
$ws.Range("H24").Value = 18408
$ws.Range("I24").Value = 112
$ws.Range("K24").Value = 112
$ws.Range("M24").Value = 58
$ws.Range("H55").Value = 28532.2
$ws.Range("I55").Value = 5998.3335
$ws.Range("J55").Value = 62333
$ws.Range("K55").Value = 5998.3335
$ws.Range("L55").Value = 62333
$ws.Range("M55").Value = -5683.3335
$ws.Range("N55").Value = -62963
$ws.Range("H93").Value = 22339.2
$ws.Range("I93").Value = 16674.125
$ws.Range("J93").Value = 44999.5
$ws.Range("K93").Value = 16674.125
$ws.Range("L93").Value = 44999.5
$ws.Range("M93").Value = -14802.125
$ws.Range("N93").Value = -48743.5
$ws.Range("H99").Value = 11350.929
$ws.Range("I99").Value = 5623.125
$ws.Range("J99").Value = 18988
$ws.Range("K99").Value = 5623.125
$ws.Range("L99").Value = 18988
$ws.Range("M99").Value = -4125.125
$ws.Range("N99").Value = -21984
$ws.Range("H122").Value = 2584.4827
$ws.Range("I122").Value = 2599.1667
$ws.Range("J122").Value = 2514
$ws.Range("K122").Value = 7797.500100000001
$ws.Range("L122").Value = 7542
$ws.Range("M122").Value = -5347.500100000001
$ws.Range("N122").Value = -12442
$ws.Range("H126").Value = 11350.929
$ws.Range("I126").Value = 5623.125
$ws.Range("J126").Value = 18988
$ws.Range("K126").Value = 16869.375
$ws.Range("L126").Value = 56964
$ws.Range("M126").Value = -14399.375
$ws.Range("N126").Value = -61904
$ws.Range("H132").Value = 43483.625
$ws.Range("I132").Value = 59858.53
$ws.Range("J132").Value = 3716
$ws.Range("K132").Value = 179575.59
$ws.Range("L132").Value = 11148
$ws.Range("M132").Value = -177045.59
$ws.Range("N132").Value = -16208
$ws.Range("H133").Value = 51349
$ws.Range("I133").Value = 15000
$ws.Range("J133").Value = 63465.332
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 63465.332
$ws.Range("M133").Value = -12470
$ws.Range("N133").Value = -68525.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 2645.0833
$ws.Range("I47").Value = 2426.4546
$ws.Range("J47").Value = 5050
$ws.Range("K47").Value = 7279.3638
$ws.Range("L47").Value = 15150
$ws.Range("M47").Value = -6848.3638
$ws.Range("N47").Value = -16012
$ws.Range("H56").Value = 27783028
$ws.Range("I56").Value = 27783028
$ws.Range("K56").Value = 27783028
$ws.Range("M56").Value = -27782498
$ws.Range("H64").Value = 6240.1816
$ws.Range("I64").Value = 25506
$ws.Range("J64").Value = 4997.2256
$ws.Range("K64").Value = 76518
$ws.Range("L64").Value = 14991.6768
$ws.Range("M64").Value = -76248
$ws.Range("N64").Value = -15531.6768
$ws.Range("H67").Value = 6240.1816
$ws.Range("I67").Value = 25506
$ws.Range("J67").Value = 4997.2256
$ws.Range("K67").Value = 76518
$ws.Range("L67").Value = 14991.6768
$ws.Range("M67").Value = -75582
$ws.Range("N67").Value = -16863.6768
$ws.Range("H113").Value = 823.6667
$ws.Range("J113").Value = 823.6667
$ws.Range("L113").Value = 2471.0001
$ws.Range("N113").Value = -6811.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5034.5557
$ws.Range("I102").Value = 6063.077
$ws.Range("K102").Value = 6063.077
$ws.Range("M102").Value = -4441.077
$ws.Range("H126").Value = 3578
$ws.Range("I126").Value = 2656
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 7968
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -5498
$ws.Range("N126").Value = -18440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1160.4445
$ws.Range("I22").Value = 596.3333
$ws.Range("K22").Value = 596.3333
$ws.Range("M22").Value = -301.3333
$ws.Range("H27").Value = 1160.4445
$ws.Range("I27").Value = 596.3333
$ws.Range("K27").Value = 596.3333
$ws.Range("M27").Value = -489.3333
$ws.Range("H46").Value = 4257.25
$ws.Range("I46").Value = 648.5
$ws.Range("J46").Value = 7866
$ws.Range("K46").Value = 648.5
$ws.Range("L46").Value = 7866
$ws.Range("M46").Value = -460.5
$ws.Range("N46").Value = -8242
$ws.Range("H122").Value = 4476.125
$ws.Range("I122").Value = 3146.625
$ws.Range("K122").Value = 9439.875
$ws.Range("M122").Value = -6989.875
$ws.Range("H132").Value = 5178.727
$ws.Range("I132").Value = 5046.6
$ws.Range("K132").Value = 15139.8
$ws.Range("M132").Value = -12609.8
$ws.Range("H136").Value = 2466.9473
$ws.Range("I136").Value = 2805.2
$ws.Range("J136").Value = 1198.5
$ws.Range("K136").Value = 8415.599999999999
$ws.Range("L136").Value = 3595.5
$ws.Range("M136").Value = -5865.599999999999
$ws.Range("N136").Value = -8695.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1410
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H113").Value = 893.069
$ws.Range("I113").Value = 671.9375
$ws.Range("K113").Value = 2015.8125
$ws.Range("M113").Value = 154.1875
$ws.Range("H122").Value = 50767.965
$ws.Range("I122").Value = 63754.094
$ws.Range("K122").Value = 191262.282
$ws.Range("M122").Value = -188812.282
$ws.Range("H132").Value = 19104.82
$ws.Range("I132").Value = 21011.514
$ws.Range("K132").Value = 63034.542
$ws.Range("M132").Value = -60504.542
$ws.Range("H136").Value = 20695.734
$ws.Range("I136").Value = 24158.965
$ws.Range("J136").Value = 4534
$ws.Range("K136").Value = 72476.895
$ws.Range("L136").Value = 13602
$ws.Range("M136").Value = -69926.895
$ws.Range("N136").Value = -18702
